$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "+" marker in column S for rows 2 and 5-16 (rows 3,4 already had it)
$ws.Range("S2").Value = "+"
$ws.Range("S5:S16").Value = "+"

# Update the view: scroll to show column E first, and select S17 as active cell
$ws.Range("S17").Select()
$excel.ActiveWindow.ScrollColumn = 5
